# Append a new data row (row 3) to the "Artfynd" sheet, mirroring the
# schema already used by the existing data row (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

# --- genuine numeric cells -------------------------------------------
$ws.Cells.Item($row, 1).Value  = 112330310      # A  Id
$ws.Cells.Item($row, 2).Value  = 56346          # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 100067         # E  TaxonId
$ws.Cells.Item($row, 17).Value = 399012         # Q  Ost
$ws.Cells.Item($row, 18).Value = 6202751        # R  Nord
$ws.Cells.Item($row, 19).Value = 25             # S  Noggrannhet

# --- plain text cells (no ambiguity with numbers/dates/times) --------
$ws.Cells.Item($row, 3).Value  = "Ovaliderad"            # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value  = "NT"                    # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Havsörn"                # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Haliaeetus albicilla"   # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Linnaeus, 1758)"       # H  Auktor
$ws.Cells.Item($row, 13).Value = "förbiflygande"          # M  Aktivitet
$ws.Cells.Item($row, 16).Value = "Bögerup hygget , Sk"    # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Skåne"                  # T  Län
$ws.Cells.Item($row, 21).Value = "Eslöv"                  # U  Kommun
$ws.Cells.Item($row, 22).Value = "Skåne"                  # V  Provins
$ws.Cells.Item($row, 23).Value = "Billinge"                # W  Församling
$ws.Cells.Item($row, 26).Value = "14:21"                   # Z  Starttid
$ws.Cells.Item($row, 28).Value = "14:21"                   # AB Sluttid
$ws.Cells.Item($row, 49).Value = "Jonny Johansson"         # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Jonny Johansson"         # AX Observatörer

# --- text cells that look numeric/date-like: force text via a
#     quote-prefixed entry, then drop the quote-prefix style again so
#     the saved cell carries plain text with no leftover formatting. --
function Set-TextValue($col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue 9  "2"            # I  Antal
Set-TextValue 25 "2023-09-26"   # Y  Startdatum
Set-TextValue 27 "2023-09-26"   # AA Slutdatum

# --- boolean cells ---------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan

# --- blank (but present) text cells -----------------------------------
# These columns carry an explicit empty-text value in the source row
# (K, L, N, AT, AY). Force text type via a quote-prefixed empty entry,
# then restore the default "Normal" style so no quote-prefix formatting
# lingers on the cell.
foreach ($col in 11, 12, 14, 46, 51) {
    Set-TextValue $col ""
}
